$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values (Sending cluster / Target cluster expression &
# derived NATMI edge metrics) as produced by the refreshed "natmiOut_TPM"
# computation. Only the numeric value cells change; labels/headers/styles
# are untouched.
$changes = @(
    @("G2", 21.83470933333334),
    @("H2", 65.50412800000001),
    @("I2", 0.1994658397831471),
    @("J2", 0.1994658397831471),
    @("M2", 9.423852333333334),
    @("N2", 28.271557),
    @("O2", 0.06654336290212845),
    @("P2", 0.06654336290212845),
    @("Q2", 205.7670764985885),
    @("R2", 1851.903688487296),
    @("S2", 0.01327312776326777),
    @("T2", 0.01327312776326777),
    @("G3", 21.83470933333334),
    @("H3", 65.50412800000001),
    @("I3", 0.1994658397831471),
    @("J3", 0.1994658397831471),
    @("O3", 0.3572423751649123),
    @("P3", 0.3572423751649123),
    @("Q3", 1104.673943924537),
    @("R3", 9942.065495320832),
    @("S3", 0.07125765036839533),
    @("T3", 0.07125765036839533),
    @("G4", 21.83470933333334),
    @("H4", 65.50412800000001),
    @("I4", 0.1994658397831471),
    @("J4", 0.1994658397831471),
    @("M4", 26.84076266666667),
    @("N4", 80.522288),
    @("O4", 0.1895270158659356),
    @("P4", 0.1895270158659356),
    @("Q4", 586.0602511116516),
    @("R4", 5274.542260004865),
    @("S4", 0.03780416538129269),
    @("T4", 0.03780416538129269),
    @("G5", 21.83470933333334),
    @("H5", 65.50412800000001),
    @("I5", 0.1994658397831471),
    @("J5", 0.1994658397831471),
    @("M5", 54.762539),
    @("N5", 164.287617),
    @("O5", 0.3866872460670236),
    @("P5", 0.3866872460670236),
    @("Q5", 1195.724121420331),
    @("R5", 10761.51709278298),
    @("S5", 0.07713089627019133),
    @("T5", 0.07713089627019132),
    @("H6", 92.79671999999999),
    @("I6", 0.282574186529459),
    @("J6", 0.282574186529459),
    @("M6", 9.423852333333334),
    @("N6", 28.271557),
    @("O6", 0.06654336290212845),
    @("P6", 0.06654336290212845),
    @("Q6", 291.5008620992267),
    @("R6", 2623.50775889304),
    @("S6", 0.01880343664100353),
    @("T6", 0.01880343664100352),
    @("H7", 92.79671999999999),
    @("I7", 0.282574186529459),
    @("J7", 0.282574186529459),
    @("O7", 0.3572423751649123),
    @("P7", 0.3572423751649123),
    @("S7", 0.1009474735560769),
    @("T7", 0.1009474735560769),
    @("H8", 92.79671999999999),
    @("I8", 0.282574186529459),
    @("J8", 0.282574186529459),
    @("M8", 26.84076266666667),
    @("N8", 80.522288),
    @("O8", 0.1895270158659356),
    @("P8", 0.1895270158659356),
    @("Q8", 830.2449125883733),
    @("R8", 7472.20421329536),
    @("S8", 0.05355544233367262),
    @("T8", 0.05355544233367262),
    @("H9", 92.79671999999999),
    @("I9", 0.282574186529459),
    @("J9", 0.282574186529459),
    @("M9", 54.762539),
    @("N9", 164.287617),
    @("O9", 0.3866872460670236),
    @("P9", 0.3866872460670236),
    @("Q9", 1693.92799935736),
    @("R9", 15245.35199421624),
    @("S9", 0.109267833998706),
    @("T9", 0.1092678339987059),
    @("G10", 47.70664233333334),
    @("H10", 143.119927),
    @("I10", 0.435812784634851),
    @("J10", 0.435812784634851),
    @("M10", 9.423852333333334),
    @("N10", 28.271557),
    @("O10", 0.06654336290212845),
    @("P10", 0.06654336290212845),
    @("Q10", 449.5803526684822),
    @("R10", 4046.22317401634),
    @("S10", 0.02900044828534404),
    @("T10", 0.02900044828534404),
    @("G11", 47.70664233333334),
    @("H11", 143.119927),
    @("I11", 0.435812784634851),
    @("J11", 0.435812784634851),
    @("O11", 0.3572423751649123),
    @("P11", 0.3572423751649123),
    @("Q11", 2413.601387278705),
    @("R11", 21722.41248550834),
    @("S11", 0.1556907943101886),
    @("T11", 0.1556907943101886),
    @("G12", 47.70664233333334),
    @("H12", 143.119927),
    @("I12", 0.435812784634851),
    @("J12", 0.435812784634851),
    @("M12", 26.84076266666667),
    @("N12", 80.522288),
    @("O12", 0.1895270158659356),
    @("P12", 0.1895270158659356),
    @("Q12", 1280.482664492553),
    @("R12", 11524.34398043298),
    @("S12", 0.08259829654806697),
    @("T12", 0.08259829654806697),
    @("G13", 47.70664233333334),
    @("H13", 143.119927),
    @("I13", 0.435812784634851),
    @("J13", 0.435812784634851),
    @("M13", 54.762539),
    @("N13", 164.287617),
    @("O13", 0.3866872460670236),
    @("P13", 0.3866872460670236),
    @("Q13", 2612.536861338218),
    @("R13", 23512.83175204396),
    @("S13", 0.1685232454912514),
    @("T13", 0.1685232454912514),
    @("G14", 8.992316666666666),
    @("H14", 26.97695),
    @("I14", 0.08214718905254291),
    @("J14", 0.08214718905254291),
    @("M14", 9.423852333333334),
    @("N14", 28.271557),
    @("O14", 0.06654336290212845),
    @("P14", 0.06654336290212845),
    @("Q14", 84.74226440123888),
    @("R14", 762.68037961115),
    @("S14", 0.005466350212513117),
    @("T14", 0.005466350212513117),
    @("G15", 8.992316666666666),
    @("H15", 26.97695),
    @("I15", 0.08214718905254291),
    @("J15", 0.08214718905254291),
    @("O15", 0.3572423751649123),
    @("P15", 0.3572423751649123),
    @("Q15", 454.9443624614777),
    @("R15", 4094.4992621533),
    @("S15", 0.02934645693025151),
    @("T15", 0.02934645693025151),
    @("G16", 8.992316666666666),
    @("H16", 26.97695),
    @("I16", 0.08214718905254291),
    @("J16", 0.08214718905254291),
    @("M16", 26.84076266666667),
    @("N16", 80.522288),
    @("O16", 0.1895270158659356),
    @("P16", 0.1895270158659356),
    @("Q16", 241.3606374735111),
    @("R16", 2172.2457372616),
    @("S16", 0.01556911160290331),
    @("T16", 0.01556911160290331),
    @("G17", 8.992316666666666),
    @("H17", 26.97695),
    @("I17", 0.08214718905254291),
    @("J17", 0.08214718905254291),
    @("M17", 54.762539),
    @("N17", 164.287617),
    @("O17", 0.3866872460670236),
    @("P17", 0.3866872460670236),
    @("Q17", 492.4420921586832),
    @("R17", 4431.97882942815),
    @("S17", 0.03176527030687497),
    @("T17", 0.03176527030687497)
)

foreach ($change in $changes) {
    $cellRef = $change[0]
    $newValue = $change[1]
    $ws.Range($cellRef).Value2 = $newValue
}

Write-Host "Updated $($changes.Count) cells."
